# Adds "Mean increase" / "Median increase" computations to the CPU processed
# results sheet for the Low/Medium/High/All frequency blocks, and renumbers
# the hidden _xlchart.v1.* defined names the way Excel does after the chart
# source ranges were touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

function Set-Header($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $true
}

# ---------------------------------------------------------------------------
# Low frequency block (rows 17-32): mean/median at row 19, headers at row 18
# ---------------------------------------------------------------------------
Set-Header "D18" "Mean increase"
Set-Header "F18" "Median increase"
$ws.Range("D19").Formula = "=((102.889693/95.321842)*100)-100"
$ws.Range("F19").Formula = "=((102.092058/95.22216)*100)-100"

# ---------------------------------------------------------------------------
# Medium frequency block (rows 49-64): mean/median at row 51, headers row 50
# ---------------------------------------------------------------------------
Set-Header "D50" "Mean increase"
Set-Header "F50" "Median increase"
$ws.Range("D51").Formula = "=((105.263347/95.321842)*100)-100"
$ws.Range("F51").Formula = "=((105.110784/95.22216)*100)-100"

# ---------------------------------------------------------------------------
# High frequency block (rows 81-96): mean/median at row 83, headers row 82
# ---------------------------------------------------------------------------
Set-Header "D82" "Mean increase"
Set-Header "F82" "Median increase"
$ws.Range("D83").Formula = "=((114.318401/95.321842)*100)-100"
$ws.Range("F83").Formula = "=((114.037416/95.22216)*100)-100"

# ---------------------------------------------------------------------------
# All block (rows 97-187): headers row 113, averages of the three increases
# at row 114
# ---------------------------------------------------------------------------
Set-Header "D113" "Mean increase"
Set-Header "F113" "Median increase"
$ws.Range("D114").Formula = "=(D19+D51+D83)/3"
$ws.Range("F114").Formula = "=(F19+F51+F83)/3"

# ---------------------------------------------------------------------------
# Renumber the hidden chart source defined names: Excel rewrote v1.2..v1.9
# in-place when the chart data ranges were re-touched (Medium <-> High <->
# All blocks swap index groups; Low/v1.0-1 stay put).
# ---------------------------------------------------------------------------
$wb.Names.Item("_xlchart.v1.2").RefersTo = "=Blad1!`$A`$34:`$A`$63"
$wb.Names.Item("_xlchart.v1.3").RefersTo = "=Blad1!`$B`$33"
$wb.Names.Item("_xlchart.v1.4").RefersTo = "=Blad1!`$B`$34:`$B`$63"
$wb.Names.Item("_xlchart.v1.5").RefersTo = "=Blad1!`$A`$98:`$A`$187"
$wb.Names.Item("_xlchart.v1.6").RefersTo = "=Blad1!`$B`$98:`$B`$187"
$wb.Names.Item("_xlchart.v1.7").RefersTo = "=Blad1!`$A`$66:`$A`$95"
$wb.Names.Item("_xlchart.v1.8").RefersTo = "=Blad1!`$B`$65"
$wb.Names.Item("_xlchart.v1.9").RefersTo = "=Blad1!`$B`$66:`$B`$95"

# ---------------------------------------------------------------------------
# Leave the cursor where the author apparently left it when saving.
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A96"), $false)
$ws.Range("I116").Select()
